# MCL questions.xlsx edit script
#
# Summary of the change being applied (see commit message: "Fixes and
# cleaned up treatment definitions."):
#   - The existing "Web Parameters" sheet (40 rows of old/legacy treatment
#     definitions) is preserved as a new sheet named "Web Parameters v1",
#     placed immediately after the original tab.
#   - The original "Web Parameters" tab is then rewritten in place with a
#     small, cleaned-up set of only 4 treatment rows that exercise the new
#     "word" view_type and populate the previously-unused graph-size columns
#     (O:R) plus derived width/height columns (S:T) and a comment column (U).

$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("Web Parameters")

# 1) Duplicate the sheet (placed right after itself) and rename the copy so
#    the old data is preserved under the new name "Web Parameters v1".
$src.Copy([System.Reflection.Missing]::Value, $src) | Out-Null
$copy = $wb.Worksheets.Item("Web Parameters (2)")
$copy.Name = "Web Parameters v1"

# 2) Wipe out the old data rows (2:40) on the original tab - only the header
#    row (row 1) is kept as-is.
$src.Range("A2:U40").Clear() | Out-Null

# 3) Write the new, cleaned-up treatment rows.
#    Columns: A=treatment_id B=position C=view_type D=titration
#             E=amount_earlier F=time_earlier G=date_earlier(unused)
#             H=amount_later I=time_later J=date_later(unused)
#             K=max_amount L=max_time M=horizontal_pixels N=vertical_pixels
#             O=left_margin_width_in P=bottom_margin_height_in
#             Q=graph_width_in R=graph_height_in
#             S=width_in (=Q+O) T=height_in (=R+P) U=comment
$rows = @(
    @{A=1; B=1; C="word";     D="none";          E=500; F=2; H=1000; I=3;  K=1000; L=10; M=480; N=480; O=0.5; P=0.5; Q=6; R=6; U="Read 2001 example, absolute size"},
    @{A=2; B=1; C="barchart"; D="none";          E=500; F=2; H=1000; I=3;  K=1000; L=10; M=480; N=480; O=0.5; P=0.5; Q=6; R=6; U="Read 2001 example, absolute size"},
    @{A=3; B=1; C="word";     D="earlierAmount"; E=500; F=2; H=1000; I=3;  K=1000; L=10; M=480; N=480; O=0.5; P=0.5; Q=6; R=6; U="Read 2001 example, absolute size"},
    @{A=4; B=1; C="barchart"; D="earlierAmount"; E=500; F=2; H=1000; I=10; K=1000; L=10; M=480; N=480; O=0.5; P=0.5; Q=6; R=6; U="Read 2001 example, absolute size"}
)

$r = 2
foreach ($row in $rows) {
    $src.Range("A$r").Value = $row.A
    $src.Range("B$r").Value = $row.B
    $src.Range("C$r").Value = $row.C
    $src.Range("D$r").Value = $row.D
    $src.Range("E$r").Value = $row.E
    $src.Range("F$r").Value = $row.F
    $src.Range("H$r").Value = $row.H
    $src.Range("I$r").Value = $row.I
    $src.Range("K$r").Value = $row.K
    $src.Range("L$r").Value = $row.L
    $src.Range("M$r").Value = $row.M
    $src.Range("N$r").Value = $row.N
    $src.Range("O$r").Value = $row.O
    $src.Range("P$r").Value = $row.P
    $src.Range("Q$r").Value = $row.Q
    $src.Range("R$r").Value = $row.R
    $src.Range("S$r").Formula = "=Q$r+O$r"
    $src.Range("T$r").Formula = "=R$r+P$r"
    $src.Range("U$r").Value = $row.U
    $r++
}

# 4) Restore the view state: the new "Web Parameters" tab stays the
#    selected/active tab (as it was before the edit), with the same frozen
#    header pane but a different remembered selected cell; "Web Parameters
#    v1" is left on its own remembered selection and is not the active tab.
#    NOTE: selecting a range on a sheet implicitly activates that sheet, so
#    the final Select()/Activate() calls must target $src last.
$copy.Range("D3").Select() | Out-Null
$src.Range("H19").Select() | Out-Null
$src.Activate() | Out-Null

Write-Host "done"
